# Auto-generated Excel COM-interop script to apply the market-price data refresh
# to the Spriggan Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 6961.6665
$ws.Range("I9").Value = 8525.583000000001
$ws.Range("K9").Value = 8525.583000000001
$ws.Range("M9").Value = -8356.583000000001
$ws.Range("H11").Value = 1007.7273
$ws.Range("I11").Value = 1007.7273
$ws.Range("K11").Value = 1007.7273
$ws.Range("M11").Value = -867.7273
$ws.Range("H12").Value = 33516
$ws.Range("I12").Value = 33516
$ws.Range("K12").Value = 33516
$ws.Range("M12").Value = -33346
$ws.Range("H17").Value = 195173.44
$ws.Range("I17").Value = 336
$ws.Range("J17").Value = 204165.94
$ws.Range("K17").Value = 1008
$ws.Range("L17").Value = 612497.8200000001
$ws.Range("M17").Value = -840
$ws.Range("N17").Value = -612833.8200000001
$ws.Range("H41").Value = 2481.3333
$ws.Range("I41").Value = 2997.25
$ws.Range("K41").Value = 2997.25
$ws.Range("M41").Value = -2557.25
$ws.Range("H62").Value = 3012.6667
$ws.Range("I62").Value = 2552
$ws.Range("K62").Value = 2552
$ws.Range("M62").Value = -1928
$ws.Range("H65").Value = 3012.6667
$ws.Range("I65").Value = 2552
$ws.Range("K65").Value = 12760
$ws.Range("M65").Value = -9640
$ws.Range("H70").Value = 5271.7896
$ws.Range("J70").Value = 7911.1113
$ws.Range("L70").Value = 23733.3339
$ws.Range("N70").Value = -24273.3339
$ws.Range("H73").Value = 5271.7896
$ws.Range("J73").Value = 7911.1113
$ws.Range("L73").Value = 23733.3339
$ws.Range("N73").Value = -25605.3339
$ws.Range("H98").Value = 1465
$ws.Range("I98").Value = 1058.5714
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 1058.5714
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = 439.4286
$ws.Range("N98").Value = -12996
$ws.Range("H99").Value = 2431.5833
$ws.Range("I99").Value = 186.77777
$ws.Range("K99").Value = 560.33331
$ws.Range("M99").Value = 937.66669
$ws.Range("H121").Value = 1199.6666
$ws.Range("J121").Value = 1199.6666
$ws.Range("L121").Value = 3598.9998
$ws.Range("N121").Value = -7092.9998
$ws.Range("H122").Value = 1465
$ws.Range("I122").Value = 1058.5714
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 3175.7142
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -725.7142000000003
$ws.Range("N122").Value = -34900
$ws.Range("H127").Value = 1556.5
$ws.Range("I127").Value = 1307.8
$ws.Range("J127").Value = 2800
$ws.Range("K127").Value = 3923.4
$ws.Range("L127").Value = 8400
$ws.Range("M127").Value = 1036.6
$ws.Range("N127").Value = -18320
$ws.Range("H129").Value = 9728.909
$ws.Range("I129").Value = 1780.1111
$ws.Range("J129").Value = 45498.5
$ws.Range("K129").Value = 5340.3333
$ws.Range("L129").Value = 136495.5
$ws.Range("M129").Value = -340.3333000000002
$ws.Range("N129").Value = -146495.5
$ws.Range("H132").Value = 3285.0625
$ws.Range("I132").Value = 3170.7334
$ws.Range("K132").Value = 9512.200199999999
$ws.Range("M132").Value = -6982.200199999999
$ws.Range("H135").Value = 26316266
$ws.Range("I135").Value = 27778268
$ws.Range("J135").Value = 226
$ws.Range("K135").Value = 250004412
$ws.Range("L135").Value = 2034
$ws.Range("M135").Value = -250001877
$ws.Range("N135").Value = -7104
$ws.Range("H137").Value = 3690.7058
$ws.Range("I137").Value = 3303.6667
$ws.Range("K137").Value = 9911.000100000001
$ws.Range("M137").Value = -7361.000100000001
$ws.Range("H138").Value = 2376.8918
$ws.Range("I138").Value = 2421.5715
$ws.Range("K138").Value = 7264.7145
$ws.Range("M138").Value = -2124.7145
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 2847
$ws.Range("I141").Value = 3084.3333
$ws.Range("J141").Value = 2135
$ws.Range("K141").Value = 9252.999899999999
$ws.Range("L141").Value = 6405
$ws.Range("M141").Value = -4072.999899999999
$ws.Range("N141").Value = -16765

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2710.205
$ws.Range("I32").Value = 2759.7715
$ws.Range("J32").Value = 2276.5
$ws.Range("K32").Value = 2759.7715
$ws.Range("L32").Value = 2276.5
$ws.Range("M32").Value = -2472.7715
$ws.Range("N32").Value = -2850.5
$ws.Range("H45").Value = 3426.82
$ws.Range("I45").Value = 3349.7954
$ws.Range("K45").Value = 3349.7954
$ws.Range("M45").Value = -2972.7954
$ws.Range("H122").Value = 4124.56
$ws.Range("I122").Value = 3254.75
$ws.Range("K122").Value = 9764.25
$ws.Range("M122").Value = -7314.25
$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -90060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1114.5714
$ws.Range("I20").Value = 1201
$ws.Range("J20").Value = 898.5
$ws.Range("K20").Value = 1201
$ws.Range("L20").Value = 898.5
$ws.Range("M20").Value = -954
$ws.Range("N20").Value = -1392.5
$ws.Range("H86").Value = 4001
$ws.Range("I86").Value = 4194.6
$ws.Range("J86").Value = 3807.4
$ws.Range("K86").Value = 4194.6
$ws.Range("L86").Value = 3807.4
$ws.Range("M86").Value = -3071.6
$ws.Range("N86").Value = -6053.4
$ws.Range("H89").Value = 4001
$ws.Range("I89").Value = 4194.6
$ws.Range("J89").Value = 3807.4
$ws.Range("K89").Value = 20973
$ws.Range("L89").Value = 19037
$ws.Range("M89").Value = -15357
$ws.Range("N89").Value = -30269
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H134").Value = 21744316
$ws.Range("I134").Value = 23814760
$ws.Range("K134").Value = 71444280
$ws.Range("M134").Value = -71441745

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 36001
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H60").Value = 4080
$ws.Range("I60").Value = 4080
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 4080
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -3569
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 36001
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H86").Value = 9596.450000000001
$ws.Range("I86").Value = 6977.5
$ws.Range("K86").Value = 6977.5
$ws.Range("M86").Value = -5854.5
$ws.Range("H89").Value = 9596.450000000001
$ws.Range("I89").Value = 6977.5
$ws.Range("K89").Value = 34887.5
$ws.Range("M89").Value = -29271.5
$ws.Range("H94").Value = 2046.8667
$ws.Range("J94").Value = 1736.8889
$ws.Range("L94").Value = 1736.8889
$ws.Range("N94").Value = -2638.8889
$ws.Range("H122").Value = 3230.3333
$ws.Range("I122").Value = 3477.4
$ws.Range("J122").Value = 1995
$ws.Range("K122").Value = 10432.2
$ws.Range("L122").Value = 5985
$ws.Range("M122").Value = -7982.200000000001
$ws.Range("N122").Value = -10885
$ws.Range("H134").Value = 25002580
$ws.Range("I134").Value = 27780422
$ws.Range("J134").Value = 2011
$ws.Range("K134").Value = 83341266
$ws.Range("L134").Value = 6033
$ws.Range("M134").Value = -83338731
$ws.Range("N134").Value = -11103

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 340.5
$ws.Range("I14").Value = 340.5
$ws.Range("K14").Value = 1021.5
$ws.Range("M14").Value = -848.5
$ws.Range("H17").Value = 2216.077
$ws.Range("I17").Value = 2278.5557
$ws.Range("K17").Value = 6835.6671
$ws.Range("M17").Value = -6666.6671
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H69").Value = 936.875
$ws.Range("J69").Value = 1166.6666
$ws.Range("L69").Value = 3499.9998
$ws.Range("N69").Value = -5121.9998
$ws.Range("H72").Value = 936.875
$ws.Range("J72").Value = 1166.6666
$ws.Range("L72").Value = 10499.9994
$ws.Range("N72").Value = -18611.9994
$ws.Range("H103").Value = 446
$ws.Range("I103").Value = 446
$ws.Range("K103").Value = 1338
$ws.Range("M103").Value = -459
$ws.Range("H122").Value = 1724.1428
$ws.Range("I122").Value = 1439.6666
$ws.Range("K122").Value = 12956.9994
$ws.Range("M122").Value = -10506.9994
$ws.Range("H129").Value = 4162.4375
$ws.Range("J129").Value = 4471.4287
$ws.Range("L129").Value = 13414.2861
$ws.Range("N129").Value = -23414.2861
$ws.Range("H139").Value = 2789.3333
$ws.Range("I139").Value = 2138
$ws.Range("J139").Value = 8000
$ws.Range("K139").Value = 6414
$ws.Range("L139").Value = 24000
$ws.Range("M139").Value = -1274
$ws.Range("N139").Value = -34280
$ws.Range("H140").Value = 3791.4
$ws.Range("I140").Value = 3989.25
$ws.Range("K140").Value = 11967.75
$ws.Range("M140").Value = -6787.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1014.5
$ws.Range("I97").Value = 1039.2354
$ws.Range("J97").Value = 930.4
$ws.Range("K97").Value = 1039.2354
$ws.Range("L97").Value = 930.4
$ws.Range("M97").Value = -543.2354
$ws.Range("N97").Value = -1922.4
$ws.Range("H122").Value = 4802.077
$ws.Range("I122").Value = 3285.5833
$ws.Range("K122").Value = 9856.749899999999
$ws.Range("M122").Value = -7406.749899999999
$ws.Range("H140").Value = 99984.5
$ws.Range("J140").Value = 99984.5
$ws.Range("L140").Value = 99984.5
$ws.Range("N140").Value = -110344.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4379.8096
$ws.Range("I7").Value = 4388.316
$ws.Range("J7").Value = 4299
$ws.Range("K7").Value = 4388.316
$ws.Range("L7").Value = 4299
$ws.Range("M7").Value = -4276.316
$ws.Range("N7").Value = -4523
$ws.Range("H16").Value = 2088.1875
$ws.Range("I16").Value = 1723.8462
$ws.Range("J16").Value = 3667
$ws.Range("K16").Value = 1723.8462
$ws.Range("L16").Value = 3667
$ws.Range("M16").Value = -1553.8462
$ws.Range("N16").Value = -4007
$ws.Range("H126").Value = 4379.8096
$ws.Range("I126").Value = 4388.316
$ws.Range("J126").Value = 4299
$ws.Range("K126").Value = 13164.948
$ws.Range("L126").Value = 12897
$ws.Range("M126").Value = -10694.948
$ws.Range("N126").Value = -17837
$ws.Range("H132").Value = 20011244
$ws.Range("I132").Value = 25275026
$ws.Range("K132").Value = 75825078
$ws.Range("M132").Value = -75822548
$ws.Range("H136").Value = 1874.95
$ws.Range("I136").Value = 1834.8667
$ws.Range("J136").Value = 1995.2
$ws.Range("K136").Value = 5504.6001
$ws.Range("L136").Value = 5985.6
$ws.Range("M136").Value = -2954.6001
$ws.Range("N136").Value = -11085.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4224.75
$ws.Range("I81").Value = 4224.75
$ws.Range("K81").Value = 8449.5
$ws.Range("M81").Value = -7388.5
$ws.Range("H84").Value = 4224.75
$ws.Range("I84").Value = 4224.75
$ws.Range("K84").Value = 42247.5
$ws.Range("M84").Value = -36943.5
$ws.Range("H96").Value = 1319.6
$ws.Range("J96").Value = 1600
$ws.Range("L96").Value = 1600
$ws.Range("N96").Value = -4346
$ws.Range("H107").Value = 803.96
$ws.Range("I107").Value = 677.8889
$ws.Range("J107").Value = 1128.1428
$ws.Range("K107").Value = 2033.6667
$ws.Range("L107").Value = 3384.4284
$ws.Range("M107").Value = -113.6667000000002
$ws.Range("N107").Value = -7224.428400000001
$ws.Range("H126").Value = 2436.3872
$ws.Range("I126").Value = 2372.4285
$ws.Range("J126").Value = 3033.3333
$ws.Range("K126").Value = 7117.2855
$ws.Range("L126").Value = 9099.999899999999
$ws.Range("M126").Value = -4647.2855
$ws.Range("N126").Value = -14039.9999
$ws.Range("H132").Value = 10002621
$ws.Range("I132").Value = 11365137
$ws.Range("J132").Value = 10838.833
$ws.Range("K132").Value = 34095411
$ws.Range("L132").Value = 32516.499
$ws.Range("M132").Value = -34092881
$ws.Range("N132").Value = -37576.499
$ws.Range("H136").Value = 38464504
$ws.Range("I136").Value = 41669720
$ws.Range("J136").Value = 1900
$ws.Range("K136").Value = 125009160
$ws.Range("L136").Value = 5700
$ws.Range("M136").Value = -125006610
$ws.Range("N136").Value = -10800
$ws.Range("H141").Value = 47993.89
$ws.Range("J141").Value = 47993.125
$ws.Range("L141").Value = 47993.125
$ws.Range("N141").Value = -58353.125
